$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update geo abbreviation values (column B) to match reverted content
$ws.Range("B6").Value = ""
$ws.Range("B11").Value = "ancsa"
$ws.Range("B12").Value = "aiaanahhl"
$ws.Range("B19").Value = "sche"
$ws.Range("B20").Value = "schs"
$ws.Range("B21").Value = "sch"
$ws.Range("B24").Value = "slupper"
$ws.Range("B25").Value = "sllower"

# Remove the geo_hierarchy value for the zcta row
$ws.Range("C23").ClearContents()

# Update the selected cell / scroll position
$ws.Range("F11").Select()
